# Fix the algorithm/conditions on filtering the status of candidates
# Insert two new candidate rows into the Active Candidates tracker:
#  - a new "Peter Abbott" (Rox / ENT AE East, 2nd Interview) row right after
#    the existing "Charles Robino" row for the same Job ID (721)
#  - a new "Tom Andrews" (Navan / ENT AE (Boston + NYC), 1st Interview) row
#    right after the "Sam Bigda-Peyton" row for Job ID 811
# This pushes all subsequent rows down and also realigns a couple of rows
# whose data had drifted out of sync with their correct row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Make room: insert one row before the current row 6, and one more
#     before the current row 15 (after the first insert shifts things down).
$ws.Rows(6).Insert()
$ws.Rows(15).Insert()

# --- Rows 1-5 are untouched. Rewrite rows 6-16 with the corrected data. ---

# Row 6: new row - Rox / ENT AE East / Peter Abbott / 2nd Interview
$ws.Range("A6").Value = 721
$ws.Range("B6").Value = "Rox"
$ws.Range("C6").Value = "ENT AE East"
$ws.Range("D6").Value = "Peter Abbott"
$ws.Range("E6").Value = "2nd Interview"
$ws.Range("F6").Value = 45973

# Row 7: Rox / RVP Sales East / TYLER FERRIER
$ws.Range("A7").Value = 725
$ws.Range("B7").Value = "Rox"
$ws.Range("C7").Value = "RVP Sales East"
$ws.Range("D7").Value = "TYLER FERRIER"
$ws.Range("E7").Value = "1st Interview"
$ws.Range("F7").Value = 45964

# Row 8: Honeycomb / VP EMEA / James Tuck
$ws.Range("A8").Value = 794
$ws.Range("B8").Value = "Honeycomb"
$ws.Range("C8").Value = "VP EMEA"
$ws.Range("D8").Value = "James Tuck"
$ws.Range("E8").Value = "1st Interview"
$ws.Range("F8").Value = 45981

# Row 9: Redwood Software / Enterprise AE US (Finance Automation) / David McDonald
$ws.Range("A9").Value = 795
$ws.Range("B9").Value = "Redwood Software"
$ws.Range("C9").Value = "Enterprise AE US (Finance Automation)"
$ws.Range("D9").Value = "David McDonald"
$ws.Range("E9").Value = "1st Interview"
$ws.Range("F9").Value = 45987

# Row 10: Redwood Software / Enterprise AE US (Finance Automation) / Tom Andrews
$ws.Range("A10").Value = 795
$ws.Range("B10").Value = "Redwood Software"
$ws.Range("C10").Value = "Enterprise AE US (Finance Automation)"
$ws.Range("D10").Value = "Tom Andrews"
$ws.Range("E10").Value = "1st Interview"
$ws.Range("F10").Value = 45986

# Row 11: Redwood Software / Enterprise AE UK (Finance Automation) / Eric Caughlin
$ws.Range("A11").Value = 796
$ws.Range("B11").Value = "Redwood Software"
$ws.Range("C11").Value = "Enterprise AE UK (Finance Automation)"
$ws.Range("D11").Value = "Eric Caughlin"
$ws.Range("E11").Value = "CV Sent"
$ws.Range("F11").Value = 45981

# Row 12: Redwood Software / Enterprise AE UK (Finance Automation) / Sam Bigda-Peyton
$ws.Range("A12").Value = 796
$ws.Range("B12").Value = "Redwood Software"
$ws.Range("C12").Value = "Enterprise AE UK (Finance Automation)"
$ws.Range("D12").Value = "Sam Bigda-Peyton"
$ws.Range("E12").Value = "1st Interview"
$ws.Range("F12").Value = 45981

# Row 13: Navan / ENT AE (Boston + NYC) / Sam Bigda-Peyton
$ws.Range("A13").Value = 811
$ws.Range("B13").Value = "Navan"
$ws.Range("C13").Value = "ENT AE (Boston + NYC)"
$ws.Range("D13").Value = "Sam Bigda-Peyton"
$ws.Range("E13").Value = "1st Interview"
$ws.Range("F13").Value = 45985

# Row 14: new row - Navan / ENT AE (Boston + NYC) / Tom Andrews
$ws.Range("A14").Value = 811
$ws.Range("B14").Value = "Navan"
$ws.Range("C14").Value = "ENT AE (Boston + NYC)"
$ws.Range("D14").Value = "Tom Andrews"
$ws.Range("E14").Value = "1st Interview"
$ws.Range("F14").Value = 45982

# Row 15: Rox / Sales Engineer (NY / Austin / SF) / Sasha Singh
$ws.Range("A15").Value = 829
$ws.Range("B15").Value = "Rox"
$ws.Range("C15").Value = "Sales Engineer (NY / Austin / SF)"
$ws.Range("D15").Value = "Sasha Singh"
$ws.Range("E15").Value = "CV Sent"
$ws.Range("F15").Value = 45987

# Row 16: Loop / Customer Operations Manager / Geoffre Lavy
$ws.Range("A16").Value = 830
$ws.Range("B16").Value = "Loop"
$ws.Range("C16").Value = "Customer Operations Manager"
$ws.Range("D16").Value = "Geoffre Lavy"
$ws.Range("E16").Value = "1st Interview"
$ws.Range("F16").Value = 45987
